# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.111.36"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.646.52"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.22"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5139"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2590"
$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.87"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07812"
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.298"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "1.657.22"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5482"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.70"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₅7771"
$ws.Range("E16").Value = "  -1.57%  "

$ws.Range("D17").Value = "26.131.50"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.25"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.473"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.990"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.090"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.903"
$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.58"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("E26").Value = "  +8.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.893"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.245"
$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04876"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.286"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.209"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.543"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.389"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9197"
$ws.Range("E35").Value = "  +2.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.604"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5505"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.119.35"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01574"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.009"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.551"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8175"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.512"
$ws.Range("E43").Value = "  -3.55%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₈125"
$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.78"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("D46").Value = "1.786.96"
$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4560"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.49"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05271"
$ws.Range("E50").Value = "  +3.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.009"
$ws.Range("E51").Value = "  +0.48%  "
